$d = $word.ActiveDocument
$d.Content.Find.Execute("BOTÃO EXCLUIR DESPESA", $true, $false, $false, $false, $false,
                         $true, 1, $false, "BOTÃO EXCLUIR", 2)
